# Monthly report: rename the "Service Status"/"Service Date" columns to
# "Visit Status"/"Visit Date" (header row 3, columns E and F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Visit Status"
$ws.Range("F3").Value = "Visit Date"

# Reproduce the reviewer's on-screen state: scrolled down one row with F6
# as the active cell.
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollRow = 2
